$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value2 = 2071.7273  # H33: was 2072.5454
$ws.Cells.Item(33, 9).Value2 = 1748.4286  # I33: was 1749.7142
$ws.Cells.Item(33, 11).Value2 = 1748.4286  # K33: was 1749.7142
$ws.Cells.Item(33, 13).Value2 = -1519.4286  # M33: was -1520.7142
$ws.Cells.Item(41, 8).Value2 = 106.181816  # H41: was 136
$ws.Cells.Item(41, 9).Value2 = 70  # I41: was 175
$ws.Cells.Item(41, 10).Value2 = 114.22222  # J41: was 126.25
$ws.Cells.Item(41, 11).Value2 = 70  # K41: was 175
$ws.Cells.Item(41, 12).Value2 = 114.22222  # L41: was 126.25
$ws.Cells.Item(41, 13).Value2 = 370  # M41: was 265
$ws.Cells.Item(41, 14).Value2 = -994.22222  # N41: was -1006.25
$ws.Cells.Item(64, 8).Value2 = 3433.3333  # H64: was 3300
$ws.Cells.Item(64, 10).Value2 = 3600  # J64: was 3000
$ws.Cells.Item(64, 12).Value2 = 3600  # L64: was 3000
$ws.Cells.Item(64, 14).Value2 = -4096  # N64: was -3496
$ws.Cells.Item(67, 8).Value2 = 3433.3333  # H67: was 3300
$ws.Cells.Item(67, 10).Value2 = 3600  # J67: was 3000
$ws.Cells.Item(67, 12).Value2 = 3600  # L67: was 3000
$ws.Cells.Item(67, 14).Value2 = -5316  # N67: was -4716
$ws.Cells.Item(75, 8).Value2 = 40725  # H75: was 33362.5
$ws.Cells.Item(75, 10).Value2 = 40725  # J75: was 33362.5
$ws.Cells.Item(75, 12).Value2 = 40725  # L75: was 33362.5
$ws.Cells.Item(75, 14).Value2 = -42597  # N75: was -35234.5
$ws.Cells.Item(78, 8).Value2 = 40725  # H78: was 33362.5
$ws.Cells.Item(78, 10).Value2 = 40725  # J78: was 33362.5
$ws.Cells.Item(78, 12).Value2 = 122175  # L78: was 100087.5
$ws.Cells.Item(78, 14).Value2 = -131535  # N78: was -109447.5
$ws.Cells.Item(98, 8).Value2 = 829.4  # H98: was 941.1539
$ws.Cells.Item(98, 9).Value2 = 312.81818  # I98: was 359.44446
$ws.Cells.Item(98, 11).Value2 = 312.81818  # K98: was 359.44446
$ws.Cells.Item(98, 13).Value2 = 1185.18182  # M98: was 1138.55554
$ws.Cells.Item(116, 8).Value2 = 10005.357  # H116: was 9511.666999999999
$ws.Cells.Item(116, 9).Value2 = 13645  # I116: was 11627.728
$ws.Cells.Item(116, 10).Value2 = 3454  # J116: was 3692.5
$ws.Cells.Item(116, 11).Value2 = 13645  # K116: was 11627.728
$ws.Cells.Item(116, 12).Value2 = 3454  # L116: was 3692.5
$ws.Cells.Item(116, 13).Value2 = -10203  # M116: was -8185.727999999999
$ws.Cells.Item(116, 14).Value2 = -10338  # N116: was -10576.5
$ws.Cells.Item(122, 8).Value2 = 829.4  # H122: was 941.1539
$ws.Cells.Item(122, 9).Value2 = 312.81818  # I122: was 359.44446
$ws.Cells.Item(122, 11).Value2 = 938.45454  # K122: was 1078.33338
$ws.Cells.Item(122, 13).Value2 = 1511.54546  # M122: was 1371.66662
$ws.Cells.Item(132, 8).Value2 = 6687.8184  # H132: was 6953.476
$ws.Cells.Item(132, 9).Value2 = 6883.3125  # I132: was 7751.643
$ws.Cells.Item(132, 10).Value2 = 6166.5  # J132: was 5357.143
$ws.Cells.Item(132, 11).Value2 = 20649.9375  # K132: was 23254.929
$ws.Cells.Item(132, 12).Value2 = 18499.5  # L132: was 16071.429
$ws.Cells.Item(132, 13).Value2 = -18119.9375  # M132: was -20724.929
$ws.Cells.Item(132, 14).Value2 = -23559.5  # N132: was -21131.429
$ws.Cells.Item(133, 8).Value2 = 59800  # H133: was 0
$ws.Cells.Item(133, 10).Value2 = 59800  # J133: was 0
$ws.Cells.Item(133, 12).Value2 = 59800  # L133: was 0
$ws.Cells.Item(133, 14).Value2 = -69920  # N133: was None
$ws.Cells.Item(138, 8).Value2 = 2024.42  # H138: was 2343.4792
$ws.Cells.Item(138, 9).Value2 = 773.26666  # I138: was 2612.25
$ws.Cells.Item(138, 10).Value2 = 2245.2117  # J138: was 2289.725
$ws.Cells.Item(138, 11).Value2 = 2319.79998  # K138: was 7836.75
$ws.Cells.Item(138, 12).Value2 = 6735.6351  # L138: was 6869.174999999999
$ws.Cells.Item(138, 13).Value2 = 2820.20002  # M138: was -2696.75
$ws.Cells.Item(138, 14).Value2 = -17015.6351  # N138: was -17149.175
$ws.Cells.Item(141, 8).Value2 = 8697.223  # H141: was 8742.714
$ws.Cells.Item(141, 9).Value2 = 3865  # I141: was 4023.75
$ws.Cells.Item(141, 10).Value2 = 18361.666  # J141: was 15034.667
$ws.Cells.Item(141, 11).Value2 = 11595  # K141: was 12071.25
$ws.Cells.Item(141, 12).Value2 = 55084.99800000001  # L141: was 45104.001
$ws.Cells.Item(141, 13).Value2 = -6415  # M141: was -6891.25
$ws.Cells.Item(141, 14).Value2 = -65444.99800000001  # N141: was -55464.001
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value2 = 393861.53  # H32: was 409054.75
$ws.Cells.Item(32, 9).Value2 = 416855.03  # I32: was 457924.9
$ws.Cells.Item(32, 10).Value2 = 35162.8  # J32: was 23523.666
$ws.Cells.Item(32, 11).Value2 = 416855.03  # K32: was 457924.9
$ws.Cells.Item(32, 12).Value2 = 35162.8  # L32: was 23523.666
$ws.Cells.Item(32, 13).Value2 = -416568.03  # M32: was -457637.9
$ws.Cells.Item(32, 14).Value2 = -35736.8  # N32: was -24097.666
$ws.Cells.Item(45, 8).Value2 = 2188.2307  # H45: was 2710.3635
$ws.Cells.Item(45, 9).Value2 = 1662.8462  # I45: was 2441.4
$ws.Cells.Item(45, 10).Value2 = 2713.6155  # J45: was 2934.5
$ws.Cells.Item(45, 11).Value2 = 1662.8462  # K45: was 2441.4
$ws.Cells.Item(45, 12).Value2 = 2713.6155  # L45: was 2934.5
$ws.Cells.Item(45, 13).Value2 = -1285.8462  # M45: was -2064.4
$ws.Cells.Item(45, 14).Value2 = -3467.6155  # N45: was -3688.5
$ws.Cells.Item(61, 8).Value2 = 2654.0322  # H61: was 2595.7188
$ws.Cells.Item(61, 9).Value2 = 1737.5  # I61: was 1687.5264
$ws.Cells.Item(61, 11).Value2 = 1737.5  # K61: was 1687.5264
$ws.Cells.Item(61, 13).Value2 = -1525.5  # M61: was -1475.5264
$ws.Cells.Item(74, 8).Value2 = 2149.7441  # H74: was 2327.8948
$ws.Cells.Item(74, 9).Value2 = 1750.6364  # I74: was 2031.4706
$ws.Cells.Item(74, 11).Value2 = 1750.6364  # K74: was 2031.4706
$ws.Cells.Item(74, 13).Value2 = -876.6364000000001  # M74: was -1157.4706
$ws.Cells.Item(77, 8).Value2 = 2149.7441  # H77: was 2327.8948
$ws.Cells.Item(77, 9).Value2 = 1750.6364  # I77: was 2031.4706
$ws.Cells.Item(77, 11).Value2 = 8753.182000000001  # K77: was 10157.353
$ws.Cells.Item(77, 13).Value2 = -4385.182000000001  # M77: was -5789.353000000001
$ws.Cells.Item(136, 8).Value2 = 2654.0322  # H136: was 2595.7188
$ws.Cells.Item(136, 9).Value2 = 1737.5  # I136: was 1687.5264
$ws.Cells.Item(136, 11).Value2 = 5212.5  # K136: was 5062.5792
$ws.Cells.Item(136, 13).Value2 = -2662.5  # M136: was -2512.5792
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value2 = 66669070  # H86: was 125003064
$ws.Cells.Item(86, 9).Value2 = 90911576  # I86: was 142859790
$ws.Cells.Item(86, 10).Value2 = 2182  # J86: was 6000
$ws.Cells.Item(86, 11).Value2 = 90911576  # K86: was 142859790
$ws.Cells.Item(86, 12).Value2 = 2182  # L86: was 6000
$ws.Cells.Item(86, 13).Value2 = -90910453  # M86: was -142858667
$ws.Cells.Item(86, 14).Value2 = -4428  # N86: was -8246
$ws.Cells.Item(89, 8).Value2 = 66669070  # H89: was 125003064
$ws.Cells.Item(89, 9).Value2 = 90911576  # I89: was 142859790
$ws.Cells.Item(89, 10).Value2 = 2182  # J89: was 6000
$ws.Cells.Item(89, 11).Value2 = 454557880  # K89: was 714298950
$ws.Cells.Item(89, 12).Value2 = 10910  # L89: was 30000
$ws.Cells.Item(89, 13).Value2 = -454552264  # M89: was -714293334
$ws.Cells.Item(89, 14).Value2 = -22142  # N89: was -41232
$ws.Cells.Item(94, 8).Value2 = 896.3333  # H94: was 926.05
$ws.Cells.Item(94, 9).Value2 = 700.7222  # I94: was 724.17645
$ws.Cells.Item(94, 11).Value2 = 700.7222  # K94: was 724.17645
$ws.Cells.Item(94, 13).Value2 = -249.7222  # M94: was -273.17645
$ws.Cells.Item(134, 8).Value2 = 1957.44  # H134: was 2044.5416
$ws.Cells.Item(134, 9).Value2 = 1660.6487  # I134: was 1775.9459
$ws.Cells.Item(134, 10).Value2 = 2802.1538  # J134: was 2948
$ws.Cells.Item(134, 11).Value2 = 4981.9461  # K134: was 5327.8377
$ws.Cells.Item(134, 12).Value2 = 8406.4614  # L134: was 8844
$ws.Cells.Item(134, 13).Value2 = -2446.9461  # M134: was -2792.8377
$ws.Cells.Item(134, 14).Value2 = -13476.4614  # N134: was -13914
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value2 = 6057.8125  # H31: was 6089.3125
$ws.Cells.Item(31, 9).Value2 = 1640.0358  # I31: was 1694.0358
$ws.Cells.Item(31, 11).Value2 = 1640.0358  # K31: was 1694.0358
$ws.Cells.Item(31, 13).Value2 = -1345.0358  # M31: was -1399.0358
$ws.Cells.Item(34, 8).Value2 = 6057.8125  # H34: was 6089.3125
$ws.Cells.Item(34, 9).Value2 = 1640.0358  # I34: was 1694.0358
$ws.Cells.Item(34, 11).Value2 = 1640.0358  # K34: was 1694.0358
$ws.Cells.Item(34, 13).Value2 = -1438.0358  # M34: was -1492.0358
$ws.Cells.Item(62, 8).Value2 = 4278.5713  # H62: was 4207.5713
$ws.Cells.Item(62, 9).Value2 = 4240  # I62: was 4377.778
$ws.Cells.Item(62, 10).Value2 = 4375  # J62: was 3901.2
$ws.Cells.Item(62, 11).Value2 = 4240  # K62: was 4377.778
$ws.Cells.Item(62, 12).Value2 = 4375  # L62: was 3901.2
$ws.Cells.Item(62, 13).Value2 = -3616  # M62: was -3753.778
$ws.Cells.Item(62, 14).Value2 = -5623  # N62: was -5149.2
$ws.Cells.Item(65, 8).Value2 = 4278.5713  # H65: was 4207.5713
$ws.Cells.Item(65, 9).Value2 = 4240  # I65: was 4377.778
$ws.Cells.Item(65, 10).Value2 = 4375  # J65: was 3901.2
$ws.Cells.Item(65, 11).Value2 = 21200  # K65: was 21888.89
$ws.Cells.Item(65, 12).Value2 = 21875  # L65: was 19506
$ws.Cells.Item(65, 13).Value2 = -18080  # M65: was -18768.89
$ws.Cells.Item(65, 14).Value2 = -28115  # N65: was -25746
$ws.Cells.Item(141, 8).Value2 = 216666.67  # H141: was 189285.72
$ws.Cells.Item(141, 10).Value2 = 220000  # J141: was 187500
$ws.Cells.Item(141, 12).Value2 = 220000  # L141: was 187500
$ws.Cells.Item(141, 14).Value2 = -230360  # N141: was -197860
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(110, 8).Value2 = 11371.559  # H110: was 11604.193
$ws.Cells.Item(110, 10).Value2 = 12051.807  # J110: was 12382.25
$ws.Cells.Item(110, 12).Value2 = 36155.421  # L110: was 37146.75
$ws.Cells.Item(110, 14).Value2 = -44335.421  # N110: was -45326.75
$ws.Cells.Item(113, 8).Value2 = 1297.2  # H113: was 1534.3636
$ws.Cells.Item(113, 10).Value2 = 1382.9231  # J113: was 1710.8889
$ws.Cells.Item(113, 12).Value2 = 4148.7693  # L113: was 5132.6667
$ws.Cells.Item(113, 14).Value2 = -8488.7693  # N113: was -9472.6667
$ws.Cells.Item(126, 8).Value2 = 3221.8  # H126: was 3038.7273
$ws.Cells.Item(126, 9).Value2 = 0  # I126: was 1715
$ws.Cells.Item(126, 10).Value2 = 3221.8  # J126: was 3332.889
$ws.Cells.Item(126, 11).Value2 = 0  # K126: was 5145
$ws.Cells.Item(126, 12).ClearContents()  # L126: was 9998.667000000001
$ws.Cells.Item(126, 13).Value2 = 9665.400000000001  # M126: was -205
$ws.Cells.Item(126, 14).Value2 = -19545.4  # N126: was -19878.667
$ws.Cells.Item(132, 8).Value2 = 2619  # H132: was 2130.5
$ws.Cells.Item(132, 9).Value2 = 2490  # I132: was 2575
$ws.Cells.Item(132, 10).Value2 = 2672.75  # J132: was 2009.2727
$ws.Cells.Item(132, 11).Value2 = 22410  # K132: was 23175
$ws.Cells.Item(132, 12).Value2 = 24054.75  # L132: was 18083.4543
$ws.Cells.Item(132, 13).Value2 = -19880  # M132: was -20645
$ws.Cells.Item(132, 14).Value2 = -29114.75  # N132: was -23143.4543
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(62, 8).Value2 = 29900  # H62: was 30666.666
$ws.Cells.Item(62, 10).Value2 = 29900  # J62: was 30666.666
$ws.Cells.Item(62, 12).Value2 = 29900  # L62: was 30666.666
$ws.Cells.Item(62, 14).Value2 = -31272  # N62: was -32038.666
$ws.Cells.Item(65, 8).Value2 = 29900  # H65: was 30666.666
$ws.Cells.Item(65, 10).Value2 = 29900  # J65: was 30666.666
$ws.Cells.Item(65, 12).Value2 = 89700  # L65: was 91999.99800000001
$ws.Cells.Item(65, 14).Value2 = -96564  # N65: was -98863.99800000001
$ws.Cells.Item(97, 8).Value2 = 1050  # H97: was 1178
$ws.Cells.Item(97, 9).Value2 = 660  # I97: was 722.5
$ws.Cells.Item(97, 11).Value2 = 660  # K97: was 722.5
$ws.Cells.Item(97, 13).Value2 = -164  # M97: was -226.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value2 = 1380.5333  # H7: was 1446.4615
$ws.Cells.Item(7, 9).Value2 = 1121.6428  # I7: was 1149.9166
$ws.Cells.Item(7, 11).Value2 = 1121.6428  # K7: was 1149.9166
$ws.Cells.Item(7, 13).Value2 = -1009.6428  # M7: was -1037.9166
$ws.Cells.Item(22, 8).Value2 = 7699.1787  # H22: was 7978.6665
$ws.Cells.Item(22, 9).Value2 = 663.4286  # I22: was 748.1667
$ws.Cells.Item(22, 10).Value2 = 10044.429  # J22: was 10044.523
$ws.Cells.Item(22, 11).Value2 = 663.4286  # K22: was 748.1667
$ws.Cells.Item(22, 12).Value2 = 10044.429  # L22: was 10044.523
$ws.Cells.Item(22, 13).Value2 = -368.4286  # M22: was -453.1667
$ws.Cells.Item(22, 14).Value2 = -10634.429  # N22: was -10634.523
$ws.Cells.Item(27, 8).Value2 = 7699.1787  # H27: was 7978.6665
$ws.Cells.Item(27, 9).Value2 = 663.4286  # I27: was 748.1667
$ws.Cells.Item(27, 10).Value2 = 10044.429  # J27: was 10044.523
$ws.Cells.Item(27, 11).Value2 = 663.4286  # K27: was 748.1667
$ws.Cells.Item(27, 12).Value2 = 10044.429  # L27: was 10044.523
$ws.Cells.Item(27, 13).Value2 = -556.4286  # M27: was -641.1667
$ws.Cells.Item(27, 14).Value2 = -10258.429  # N27: was -10258.523
$ws.Cells.Item(34, 8).Value2 = 0  # H34: was 10000
$ws.Cells.Item(34, 9).Value2 = 0  # I34: was 10000
$ws.Cells.Item(34, 11).Value2 = 0  # K34: was 10000
$ws.Cells.Item(34, 13).ClearContents()  # M34: was -9828
$ws.Cells.Item(126, 8).Value2 = 1380.5333  # H126: was 1446.4615
$ws.Cells.Item(126, 9).Value2 = 1121.6428  # I126: was 1149.9166
$ws.Cells.Item(126, 11).Value2 = 3364.9284  # K126: was 3449.7498
$ws.Cells.Item(126, 13).Value2 = -894.9284000000002  # M126: was -979.7498000000001
$ws.Cells.Item(132, 8).Value2 = 3833.875  # H132: was 4004.5652
$ws.Cells.Item(132, 9).Value2 = 3186.5  # I132: was 3438.6924
$ws.Cells.Item(132, 11).Value2 = 9559.5  # K132: was 10316.0772
$ws.Cells.Item(132, 13).Value2 = -7029.5  # M132: was -7786.0772
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value2 = 10842.25  # H45: was 10187.6
$ws.Cells.Item(45, 9).Value2 = 6569  # I45: was 7069
$ws.Cells.Item(45, 11).Value2 = 6569  # K45: was 7069
$ws.Cells.Item(45, 13).Value2 = -6078  # M45: was -6578
$ws.Cells.Item(76, 8).Value2 = 66666.664  # H76: was 68333.336
$ws.Cells.Item(76, 10).Value2 = 66666.664  # J76: was 68333.336
$ws.Cells.Item(76, 12).Value2 = 66666.664  # L76: was 68333.336
$ws.Cells.Item(76, 14).Value2 = -67296.664  # N76: was -68963.336
$ws.Cells.Item(79, 8).Value2 = 66666.664  # H79: was 68333.336
$ws.Cells.Item(79, 10).Value2 = 66666.664  # J79: was 68333.336
$ws.Cells.Item(79, 12).Value2 = 66666.664  # L79: was 68333.336
$ws.Cells.Item(79, 14).Value2 = -68850.664  # N79: was -70517.336
$ws.Cells.Item(132, 8).Value2 = 3270124.8  # H132: was 3403561.5
$ws.Cells.Item(132, 9).Value2 = 2608.5264  # I132: was 2627.0527
$ws.Cells.Item(132, 10).Value2 = 5210212.5  # J132: was 5557486.5
$ws.Cells.Item(132, 11).Value2 = 7825.5792  # K132: was 7881.158100000001
$ws.Cells.Item(132, 12).Value2 = 15630637.5  # L132: was 16672459.5
$ws.Cells.Item(132, 13).Value2 = -5295.5792  # M132: was -5351.158100000001
$ws.Cells.Item(132, 14).Value2 = -15635697.5  # N132: was -16677519.5
$ws.Cells.Item(136, 8).Value2 = 4041.1052  # H136: was 4755.143
$ws.Cells.Item(136, 9).Value2 = 3982.7693  # I136: was 4506.5454
$ws.Cells.Item(136, 10).Value2 = 4167.5  # J136: was 5666.6665
$ws.Cells.Item(136, 11).Value2 = 11948.3079  # K136: was 13519.6362
$ws.Cells.Item(136, 12).Value2 = 12502.5  # L136: was 16999.9995
$ws.Cells.Item(136, 13).Value2 = -9398.3079  # M136: was -10969.6362
$ws.Cells.Item(136, 14).Value2 = -17602.5  # N136: was -22099.9995
